$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second listing (row 3) entirely - cells shift up, shared
# strings that become unused are dropped automatically by the engine.
$ws.Range("A3:F3").EntireRow.Delete()

# The hyperlink collection keeps a stale entry for the deleted row, and
# there is no reliable in-place "retarget" operation, so rebuild the
# hyperlink for the surviving row from scratch.
$ws.Hyperlinks.Delete()

# Row 2 (the only remaining listing) - update to the new property.
$ws.Range("A2").Value = "Особняк на Александрова, 18"

$ws.Hyperlinks.Add($ws.Range("B2"), "https://osobnyaki.com//na-nikoloyamskoy-49s2")
# Re-apply the workbook's existing Hyperlink cell style so we don't end up
# with a duplicate, functionally-identical style entry.
$ws.Range("B2").Style = "Hyperlink"

# Date/Square/Price are stored as plain text in this sheet (not numbers or
# dates) - use a leading apostrophe to stop Excel from auto-converting
# them, then drop back to the default cell style.
$ws.Range("C2").Value = "'2021-12-16"
$ws.Range("D2").Value = "'1788"
$ws.Range("E2").Value = "'60067867890"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"

# Price per square meter is a genuine number.
$ws.Range("F2").Value = 54152
